$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 318 ("調子どう？" post), which shifts all subsequent
# rows (319-323) up by one and shrinks the used range from A1:C323 to A1:C322.
$ws.Rows.Item(318).Delete()
